$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Task ID 2, Registration Page): mark as Completed with 100% progress,
# and highlight the row yellow.
$ws.Range("A3:G3").Interior.Color = 65535
$ws.Range("F3").Value = "Completed"
$ws.Range("G3").Value = 100

# Row 9 (Task ID 8, Edit Profile Page): mark as Completed with 100% progress,
# and highlight the row with the theme accent6 (green) color.
$ws.Range("A9:G9").Interior.Color = 65535
$ws.Range("A9:G9").Interior.ThemeColor = 10
$ws.Range("F9").Value = "Completed"
$ws.Range("G9").Value = 100

# Update the saved view: scroll back to the top and select A9:G9.
$ws.Range("A9:G9").Select()
